# Updates cryptos list figures (price + 1h volume change) per the
# "Updated cryptos list ... with GitHub Actions" commit, including the
# ARBITRUM / HuobiToken row swap at rows 33-34.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.847.40'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '1.811.54'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.08%  '
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4659'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3696'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07352'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8688'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.38'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.27%  '
$ws.Range('D12').Value = '1.875.43'
$ws.Range('E12').Value = '  +3.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.343'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07071'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('E15').Value = '  -2.44%  '
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.003'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E18').Value = '  -0.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.87%  '
$ws.Range('D21').Value = '26.886.70'
$ws.Range('E21').Value = '  -1.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.331'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.55'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.35%  '
$ws.Range('D24').Value = '2.092.99'
$ws.Range('E24').Value = '  +2.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.896'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.107'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.287'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('E30').Value = '  -1.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08917'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7554'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.00%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.930'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.149'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.81%  '
$ws.Range('E35').Value = '  -1.59%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.095'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.77%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01951'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05253'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5330'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.906'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.167'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.59%  '
$ws.Range('E43').Value = '  -1.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1660'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.412'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4929'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.34'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.95'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06272'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.79%  '
